$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'42.587.15"
$ws.Range("E2").Value = '  -7.54%  '
$ws.Range("D3").Value = "'2.551.66"
$ws.Range("E3").Value = '  -1.96%  '
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = '  -0.17%  '
$ws.Range("D5").Value = "'296.04"
$ws.Range("E5").Value = '  -4.99%  '
$ws.Range("D6").Value = "'93.08"
$ws.Range("E6").Value = '  -6.12%  '
$ws.Range("E7").Value = '  -4.18%  '
$ws.Range("E8").Value = '  -0.04%  '
$ws.Range("D9").Value = "'0.545"
$ws.Range("E9").Value = '  -6.05%  '
$ws.Range("D10").Value = "'35.41"
$ws.Range("E10").Value = '  -9.32%  '
$ws.Range("D11").Value = "'0.0802"
$ws.Range("E11").Value = '  -4.35%  '
$ws.Range("D12").Value = "'7.67"
$ws.Range("E12").Value = '  -5.87%  '
$ws.Range("D13").Value = "'2.938.19"
$ws.Range("E13").Value = '  -2.23%  '
$ws.Range("E14").Value = '  +0.11%  '
$ws.Range("D15").Value = "'2.550.62"
$ws.Range("E15").Value = '  -1.81%  '
$ws.Range("D16").Value = "'0.863"
$ws.Range("E16").Value = '  -5.81%  '
$ws.Range("D17").Value = "'14.06"
$ws.Range("E17").Value = '  -5.44%  '
$ws.Range("D18").Value = "'42.607.66"
$ws.Range("E18").Value = '  -7.74%  '
$ws.Range("D19").Value = "'12.67"
$ws.Range("E19").Value = '  -1.61%  '
$ws.Range("D20").Value = "'6.58"
$ws.Range("E20").Value = '  -2.22%  '
$ws.Range("E21").Value = '  -4.98%  '
$ws.Range("D22").Value = "'72.09"
$ws.Range("E22").Value = '  +0.25%  '
$ws.Range("D23").Value = "'257.28"
$ws.Range("E23").Value = '  -6.50%  '
$ws.Range("D24").Value = "'2.88"
$ws.Range("E24").Value = '  -6.32%  '
$ws.Range("D25").Value = "'29.45"
$ws.Range("E25").Value = '  -1.76%  '
$ws.Range("E26").Value = '  -4.44%  '
$ws.Range("E27").Value = '  +0.26%  '
$ws.Range("E28").Value = '  -7.77%  '
$ws.Range("E29").Value = '  -3.35%  '
$ws.Range("E30").Value = '  -5.83%  '
$ws.Range("E31").Value = '  -5.85%  '
$ws.Range("B32").Value = 'Monero'
$ws.Range("C32").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D32").Value = "'150.50"
$ws.Range("E32").Value = '  -3.16%  '
$ws.Range("B33").Value = 'LidoDAOToken'
$ws.Range("C33").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D33").Value = "'3.39"
$ws.Range("E33").Value = '  -6.33%  '
$ws.Range("E34").Value = '  -3.66%  '
$ws.Range("D35").Value = "'2.73"
$ws.Range("E35").Value = '  -3.04%  '
$ws.Range("D36").Value = "'0.0790"
$ws.Range("E36").Value = '  -5.54%  '
$ws.Range("D37").Value = "'0.113"
$ws.Range("E37").Value = '  -8.25%  '
$ws.Range("D38").Value = "'24.21"
$ws.Range("E38").Value = '  +3.70%  '
$ws.Range("E39").Value = '  -3.49%  '
$ws.Range("D40").Value = "'15.73"
$ws.Range("E40").Value = '  -0.82%  '
$ws.Range("D41").Value = "'3.40"
$ws.Range("E41").Value = '  -5.85%  '
$ws.Range("D42").Value = "'0.0307"
$ws.Range("E42").Value = '  -7.12%  '
$ws.Range("B43").Value = 'RenderToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D43").Value = "'3.79"
$ws.Range("E43").Value = '  -4.49%  '
$ws.Range("B44").Value = 'Maker'
$ws.Range("C44").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D44").Value = "'2.054.62"
$ws.Range("E44").Value = '  -2.23%  '
$ws.Range("E45").Value = '  -0.12%  '
$ws.Range("D46").Value = "'84.17"
$ws.Range("E46").Value = '  -11.76%  '
$ws.Range("D47").Value = "'1.59"
$ws.Range("E47").Value = '  +3.15%  '
$ws.Range("D48").Value = "'2.796.28"
$ws.Range("E48").Value = '  -2.26%  '
$ws.Range("D49").Value = "'8.72"
$ws.Range("E49").Value = '  -9.38%  '
$ws.Range("B50").Value = 'Aave'
$ws.Range("C50").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D50").Value = "'103.13"
$ws.Range("E50").Value = '  -5.41%  '
$ws.Range("B51").Value = 'Stacks'
$ws.Range("C51").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D51").Value = "'1.67"
$ws.Range("E51").Value = '  -4.26%  '
